$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, pushing the existing rows 87-99 down to
# 88-100 (this also grows the sheet's used range from R99 to R100).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new weekly price record.
$ws.Range("A87").Value = 5
$ws.Range("B87").Value = "Macroferia Regional de Talca"
$ws.Range("C87").Value = "Maule"
$ws.Range("D87").Value = 44491
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112031
$ws.Range("G87").Value = "Poroto verde"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 100
$ws.Range("K87").Value = 42000
$ws.Range("L87").Value = 42000
$ws.Range("M87").Value = 42000
$ws.Range("N87").Value = "`$/saco 25 kilos"
$ws.Range("O87").Value = "Región del Maule"
$ws.Range("P87").Value = 1680
$ws.Range("Q87").Value = 25
$ws.Range("R87").Value = "Hortaliza"
